$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.994.53"
$ws.Range("E2").Value = "  +1.00%  "
$ws.Range("D3").Value = "1.885.83"
$ws.Range("E3").Value = "  +0.61%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "331.40"
$ws.Range("E5").Value = "  -2.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4595"
$ws.Range("E7").Value = "  -2.40%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4085"
$ws.Range("E8").Value = "  +2.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.22"
$ws.Range("E9").Value = "  -1.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07984"
$ws.Range("E10").Value = "  -0.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9898"
$ws.Range("E11").Value = "  -1.56%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.69"
$ws.Range("E12").Value = "  -2.18%  "
$ws.Range("D13").Value = "1.865.15"
$ws.Range("E13").Value = "  +0.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.900"
$ws.Range("E14").Value = "  -2.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.062"
$ws.Range("E15").Value = "  -3.20%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  +0.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.77"
$ws.Range("E17").Value = "  -1.94%  "
$ws.Range("E18").Value = "  -1.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06562"
$ws.Range("E19").Value = "  -0.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.42"
$ws.Range("E20").Value = "  -1.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9984"
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("D22").Value = "29.021.77"
$ws.Range("E22").Value = "  +1.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.400"
$ws.Range("E23").Value = "  -2.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.25"
$ws.Range("E24").Value = "  +1.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.212"
$ws.Range("E25").Value = "  -2.07%  "
$ws.Range("D26").Value = "2.120.03"
$ws.Range("E26").Value = "  +1.61%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.95"
$ws.Range("E27").Value = "  -2.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.60"
$ws.Range("E28").Value = "  -1.06%  "
$ws.Range("E29").Value = "  -1.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.408"
$ws.Range("E30").Value = "  -1.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.72"
$ws.Range("E31").Value = "  -1.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9745"
$ws.Range("E32").Value = "  -0.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09332"
$ws.Range("E33").Value = "  -2.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.600"
$ws.Range("E34").Value = "  -2.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.405"
$ws.Range("E35").Value = "  +0.94%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.274"
$ws.Range("E36").Value = "  -1.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06039"
$ws.Range("E37").Value = "  -2.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02220"
$ws.Range("E38").Value = "  -1.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.243"
$ws.Range("E39").Value = "  -2.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.178"
$ws.Range("E40").Value = "  -0.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9996"
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5758"
$ws.Range("E42").Value = "  -3.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1818"
$ws.Range("E43").Value = "  -3.64%  "
$ws.Range("E44").Value = "  -2.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.258"
$ws.Range("E45").Value = "  -0.13%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.93"
$ws.Range("E46").Value = "  -2.56%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.252"
$ws.Range("E47").Value = "  +8.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5454"
$ws.Range("E48").Value = "  -2.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.895"
$ws.Range("E49").Value = "  -3.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07024"
$ws.Range("E50").Value = "  -4.93%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "45.48"
$ws.Range("E51").Value = "  +13.35%  "
